# Realestate Update resale numbers 2025-02-12 22:35
# Appends a new data row (row 66) to the CityResaleNum sheet with the
# latest resale numbers snapshot, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66

# Columns A-D hold text values (date/time/weekday/week-number strings).
# A leading apostrophe forces Excel to treat the entry as text instead of
# auto-converting it to a date/time/number; re-applying the "Normal" style
# afterwards clears the quote-prefix formatting so the cell ends up with no
# explicit style, matching the rest of the data rows.
function Set-TextCell($r, $c, $text) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $row 1 "2025-02-12"
Set-TextCell $row 2 "22:35:53"
Set-TextCell $row 3 "Wednesday"
Set-TextCell $row 4 "06"

# Columns E-T hold numeric values.
$ws.Cells.Item($row, 5).Value = 128209
$ws.Cells.Item($row, 6).Value = 142137
$ws.Cells.Item($row, 7).Value = 169742
$ws.Cells.Item($row, 8).Value = 158907
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 144678
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191789
$ws.Cells.Item($row, 14).Value = 115091
$ws.Cells.Item($row, 15).Value = 45018
$ws.Cells.Item($row, 16).Value = 28600
$ws.Cells.Item($row, 17).Value = 65544
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 44762
$ws.Cells.Item($row, 20).Value = -1
